$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("B1").Value = "FFR_LF"
$ws.Range("C1").Value = "FFR_A"

# Copy the existing header style (bold, bordered, centered) onto the new D1 cell
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D1").Value = "A_C"

# Update row 2 values
$ws.Range("A2").Value = "params"
$ws.Range("B2").Value = 0.0225387271325573797231545114527762052603065967559814453125
$ws.Range("C2").Value = -0.6064349805449229524612064778921194374561309814453125
$ws.Range("D2").Value = 0.11774091574494739698142353745424770750105381011962890625

# Update row 3 values
$ws.Range("A3").Value = "pvalue"
$ws.Range("B3").Value = 0.00000343655427847977492487071930760134108595593716017901897430419921875
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
